$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function Replace-ParagraphXml($paraIndex, $innerXml) {
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $p.Range
    $xml = "<w:p $wns>" + $innerXml + "</w:p>"
    [void]$r.InsertXML($xml)
}

# Paragraph 1: main H1 title
Replace-ParagraphXml 1 "<w:pPr><w:pStyle w:val='Heading1'/></w:pPr><w:r><w:t>Play Dolphin Treasure for Free - Exciting Gameplay &amp; Big Payouts</w:t></w:r>"

# "What we like" bullets (paragraphs 44-47)
Replace-ParagraphXml 44 "<w:pPr><w:pStyle w:val='ListBullet'/><w:spacing w:line='240' w:lineRule='auto'/><w:ind w:left='720'/></w:pPr><w:r/><w:r><w:t>Exciting gameplay features</w:t></w:r>"
Replace-ParagraphXml 45 "<w:pPr><w:pStyle w:val='ListBullet'/><w:spacing w:line='240' w:lineRule='auto'/><w:ind w:left='720'/></w:pPr><w:r/><w:r><w:t>Vibrant and colorful graphics</w:t></w:r>"
Replace-ParagraphXml 46 "<w:pPr><w:pStyle w:val='ListBullet'/><w:spacing w:line='240' w:lineRule='auto'/><w:ind w:left='720'/></w:pPr><w:r/><w:r><w:t>Potential for significant payouts</w:t></w:r>"
Replace-ParagraphXml 47 "<w:pPr><w:pStyle w:val='ListBullet'/><w:spacing w:line='240' w:lineRule='auto'/><w:ind w:left='720'/></w:pPr><w:r/><w:r><w:t>Free spins bonus round</w:t></w:r>"

# "What we don't like" bullets (paragraphs 49-50)
Replace-ParagraphXml 49 "<w:pPr><w:pStyle w:val='ListBullet'/><w:spacing w:line='240' w:lineRule='auto'/><w:ind w:left='720'/></w:pPr><w:r/><w:r><w:t>Limited betting options</w:t></w:r>"
Replace-ParagraphXml 50 "<w:pPr><w:pStyle w:val='ListBullet'/><w:spacing w:line='240' w:lineRule='auto'/><w:ind w:left='720'/></w:pPr><w:r/><w:r><w:t>No progressive jackpot</w:t></w:r>"

# Paragraph 51: bold repeated title
Replace-ParagraphXml 51 "<w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Dolphin Treasure for Free - Exciting Gameplay &amp; Big Payouts</w:t></w:r>"

# Paragraph 52: italic meta description
Replace-ParagraphXml 52 "<w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Dolphin Treasure, an entertaining online slot game with exciting gameplay features and the potential for big payouts. Play for free now.</w:t></w:r>"

Write-Output "done"
